# Se eliminan duplicados en Veterinarios.xlsx
# Replace the two duplicate veterinarian records with the de-duplicated data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Joaquin Sosa Estrada
$ws.Range("A2").Value = "Joaquin"
$ws.Range("B2").Value = "Sosa"
$ws.Range("C2").Value = "Estrada"
$ws.Range("D2").Value = "sosa-je@gmail.com"
$ws.Range("F2").Value = 56781234

# Row 3: Matha Pastrana Rios
$ws.Range("A3").Value = "Matha"
$ws.Range("B3").Value = "Pastrana"
$ws.Range("C3").Value = "Rios"
$ws.Range("D3").Value = "rios-mpast@gmail.com"
$ws.Range("F3").Value = 21876543

# Update the active cell / selection left after the edit.
$ws.Range("E6").Select() | Out-Null
